$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.84105110168457
$ws.Range("B1").Value = 2.043452262878418
$ws.Range("C1").Value = 2.444795846939087
$ws.Range("D1").Value = 3.761769533157349
$ws.Range("E1").Value = 1.178920865058899
